$d = $word.ActiveDocument

# Metric highlight color (hex 2C3E50) expressed as a Word BGR-packed long
$metricColor = 5258796

function Highlight-Metric($paraIndex, $metricText) {
    $p = $d.Paragraphs($paraIndex)
    $searchRange = $p.Range
    $searchRange.Find.ClearFormatting()
    $found = $searchRange.Find.Execute($metricText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $searchRange.Font.Bold = 1
        $searchRange.Font.Color = $metricColor
    }
}

# Paragraph 10: "... improving demographic classification accuracy from 23% to 64%"
Highlight-Metric 10 "23%"
Highlight-Metric 10 "64%"

# Paragraph 12: "... margin of error from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87% ..."
Highlight-Metric 12 "±4.2%"
Highlight-Metric 12 "±2.1%"
Highlight-Metric 12 "71%"
Highlight-Metric 12 "87%"

# Paragraph 13: "... reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M ..."
Highlight-Metric 13 "73.5%"
Highlight-Metric 13 "$4.7M"

# Paragraph 14: "... valued over $2 trillion"
Highlight-Metric 14 "$2"

# Paragraph 34: "... reducing processing time by 57%"
Highlight-Metric 34 "57%"

# Paragraph 50: "... margin of error from ±4.2% to ±2.1%"
Highlight-Metric 50 "±4.2%"
Highlight-Metric 50 "±2.1%"

# Paragraph 51: "• Increased voter turnout prediction accuracy from 71% to 87%"
Highlight-Metric 51 "71%"
Highlight-Metric 51 "87%"

# Paragraph 52: "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"
Highlight-Metric 52 "34%"
Highlight-Metric 52 "28%"
